$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Finance")
$ws2 = $wb.Worksheets.Item("Quotation")

# ---------------------------------------------------------------------------
# Sheet "Finance" — new quotation line for "Box" + refreshed totals
# ---------------------------------------------------------------------------

# Row 72: the stray "-" text in D72 is replaced by the numeric value that used
# to live in E72, and E72 itself goes away (columns collapse by one).
$ws1.Range("D72").Value2 = 100
$ws1.Range("E72").Clear()

# Row 73: "Sheet cutting" line gets its real cost (150) entered in both the
# cost column and the "paid" column, plus a per-head formula in G.
$ws1.Range("C73").Value2 = 150
$ws1.Range("C73").Copy()
$ws1.Range("E73").PasteSpecial(-4122)   # xlPasteFormats — match C73's look
$ws1.Range("E73").Value2 = 150
$ws1.Range("G73").Formula = "=C73/3"

# Row 75: totals row — add the missing per-head total in G (style borrowed
# from the equivalent G61 total cell above).
$ws1.Range("G61").Copy()
$ws1.Range("G75").PasteSpecial(-4122)
$ws1.Range("G75").Formula = "=SUM(G72:G73)"

# Row 79: Husein's running balance before this batch.
$ws1.Range("C79").Value2 = 4620

# ---------------------------------------------------------------------------
# Sheet "Quotation" — add the "Box" line item and push the grand total down
# ---------------------------------------------------------------------------

$ws2.Rows("19:20").Insert()

$ws2.Range("A18:C18").Copy()
$ws2.Range("A19:C19").PasteSpecial(-4122)   # xlPasteFormats — match row 18
$ws2.Range("A19").Value2 = "Box"
$ws2.Range("B19").Value2 = 1
$ws2.Range("C19").Value2 = 250

$ws2.Range("A20:C20").Clear()

$ws2.Range("C21").Formula = "=SUM(C3:C19)"

# ---------------------------------------------------------------------------
# Window state — Finance becomes the active/visible tab with a cursor near
# the newly edited rows; Quotation's own selection moves off the old total.
# ---------------------------------------------------------------------------

$ws2.Range("G14").Select()
$ws1.Activate()
$ws1.Range("K72").Select()
